# Fjernet ubrukt kode i aInntekt samt ubrukte kommentarer.
# Clears the now-unused debitor_ident / Sak_Nr / FNR / debitor_ident(E) values
# that were previously hard-coded in rows 5-25 of Ark1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Clear columns A, B, C and E for rows 5 through 25 (column D was already empty).
$ws.Range("A5:C25").ClearContents()
$ws.Range("E5:E25").ClearContents()

# Update the active selection left behind in the sheet view.
$ws.Range("H14").Select()
